$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diva")

$ws.Range("A3").Value = 467
$ws.Range("A4").Value = 468
